$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "27.713.04"
Set-TextValue "E2" "  +0.42%  "
Set-TextValue "D3" "1.852.68"
Set-TextValue "E3" "  +0.57%  "
Set-TextValue "E4" "  +0.25%  "
Set-TextValue "D5" "312.43"
Set-TextValue "E5" "  -0.59%  "
Set-TextValue "E6" "  +0.28%  "
Set-TextValue "D7" "0.4280"
Set-TextValue "E7" "  +1.19%  "
Set-TextValue "D8" "0.3593"
Set-TextValue "E8" "  -1.18%  "
Set-TextValue "D9" "0.07294"
Set-TextValue "E9" "  +0.19%  "
Set-TextValue "D10" "0.8746"
Set-TextValue "E10" "  -1.48%  "
Set-TextValue "D11" "20.78"
Set-TextValue "E11" "  +0.38%  "
Set-TextValue "D12" "1.832.03"
Set-TextValue "E12" "  -0.81%  "
Set-TextValue "D13" "6.550"
Set-TextValue "E13" "  -0.13%  "
Set-TextValue "D14" "5.336"
Set-TextValue "E14" "  +0.04%  "
Set-TextValue "D15" "0.07006"
Set-TextValue "E15" "  +1.72%  "
Set-TextValue "D16" "1.005"
Set-TextValue "E16" "  +0.32%  "
Set-TextValue "D17" "79.61"
Set-TextValue "E17" "  +0.82%  "
Set-TextValue "D18" "0.000008948"
Set-TextValue "E18" "  +0.76%  "
Set-TextValue "D19" "1.003"
Set-TextValue "E19" "  +0.26%  "
Set-TextValue "D20" "15.30"
Set-TextValue "E20" "  -0.90%  "
Set-TextValue "D21" "27.729.56"
Set-TextValue "E21" "  +0.49%  "
Set-TextValue "D22" "4.997"
Set-TextValue "E22" "  +0.17%  "
Set-TextValue "D23" "10.40"
Set-TextValue "E23" "  -1.73%  "
Set-TextValue "D24" "2.048.91"
Set-TextValue "E24" "  -0.92%  "
Set-TextValue "E25" "  +4.64%  "
Set-TextValue "D26" "154.91"
Set-TextValue "E26" "  +0.76%  "
Set-TextValue "E27" "  -2.38%  "
Set-TextValue "D28" "120.46"
Set-TextValue "E28" "  -2.09%  "
Set-TextValue "D29" "5.260"
Set-TextValue "E29" "  -0.74%  "
Set-TextValue "E30" "  -0.37%  "
Set-TextValue "E31" "  -0.01%  "
Set-TextValue "D32" "0.7587"
Set-TextValue "E32" "  -1.53%  "
Set-TextValue "D33" "2.970"
Set-TextValue "E33" "  +1.92%  "
Set-TextValue "D34" "4.516"
Set-TextValue "E34" "  -1.19%  "
Set-TextValue "D35" "1.124"
Set-TextValue "E35" "  +2.50%  "
Set-TextValue "D36" "1.002"
Set-TextValue "E36" "  +0.29%  "
Set-TextValue "B37" "Hedera"
Set-TextValue "C37" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D37" "0.05439"
Set-TextValue "E37" "  +1.05%  "
Set-TextValue "B38" "TrustWalletToken"
Set-TextValue "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.108"
Set-TextValue "E38" "  +0.59%  "
Set-TextValue "E39" "  -0.48%  "
Set-TextValue "E40" "  +0.24%  "
Set-TextValue "D41" "0.1673"
Set-TextValue "E41" "  +1.23%  "
Set-TextValue "D42" "0.5097"
Set-TextValue "E42" "  +0.15%  "
Set-TextValue "E43" "  -3.59%  "
Set-TextValue "D44" "8.420"
Set-TextValue "E44" "  +1.91%  "
Set-TextValue "D45" "106.16"
Set-TextValue "E45" "  +1.72%  "
Set-TextValue "D46" "0.06527"
Set-TextValue "E46" "  -1.05%  "
Set-TextValue "D47" "10.35"
Set-TextValue "E47" "  -0.91%  "
Set-TextValue "D48" "0.4677"
Set-TextValue "E48" "  -0.93%  "
Set-TextValue "D49" "1.002"
Set-TextValue "E49" "  +0.29%  "
Set-TextValue "D50" "1.623"
Set-TextValue "E50" "  -0.52%  "
Set-TextValue "E51" "  +3.27%  "
